## Fruta / hortaliza, semanal
## Inserts one new weekly price-report row for "Mandarina" (variety Murcott)
## above the existing row 315, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 315 (rows 315..385 shift down to 316..386).
$ws.Rows.Item(315).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(315, 1).Value  = 5
$ws.Cells.Item(315, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(315, 3).Value  = "Maule"
$ws.Cells.Item(315, 4).Value  = 44889
$ws.Cells.Item(315, 5).Value  = 7
$ws.Cells.Item(315, 6).Value  = "Fruta"
$ws.Cells.Item(315, 7).Value  = 100102
$ws.Cells.Item(315, 8).Value  = "Cítricos"
$ws.Cells.Item(315, 9).Value  = 100102004
$ws.Cells.Item(315, 10).Value = "Mandarina"
$ws.Cells.Item(315, 11).Value = "Murcott"
$ws.Cells.Item(315, 12).Value = "Primera"
$ws.Cells.Item(315, 13).Value = 550
$ws.Cells.Item(315, 14).Value = 6500
$ws.Cells.Item(315, 15).Value = 7000
$ws.Cells.Item(315, 16).Value = 6727
$ws.Cells.Item(315, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(315, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(315, 19).Value = 374
$ws.Cells.Item(315, 20).Value = 18
